$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'42.759.48"
$ws.Cells.Item(2, 5).Value = "  +0.20%  "

$ws.Cells.Item(3, 4).Value = "'2.310.42"
$ws.Cells.Item(3, 5).Value = "  +0.71%  "

$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  -0.04%  "

$ws.Cells.Item(5, 4).Value = "'301.37"
$ws.Cells.Item(5, 5).Value = "  -0.17%  "

$ws.Cells.Item(6, 4).Value = "'95.25"
$ws.Cells.Item(6, 5).Value = "  -0.72%  "

$ws.Cells.Item(7, 4).Value = "'0.504"
$ws.Cells.Item(7, 5).Value = "  +0.09%  "

$ws.Cells.Item(8, 5).Value = "  +0.04%  "

$ws.Cells.Item(9, 5).Value = "  -0.96%  "

$ws.Cells.Item(10, 4).Value = "'34.13"
$ws.Cells.Item(10, 5).Value = "  -1.73%  "

$ws.Cells.Item(11, 4).Value = "'18.97"
$ws.Cells.Item(11, 5).Value = "  +2.28%  "

$ws.Cells.Item(12, 5).Value = "  +0.21%  "

$ws.Cells.Item(13, 5).Value = "  +0.01%  "

$ws.Cells.Item(14, 4).Value = "'6.72"
$ws.Cells.Item(14, 5).Value = "  -1.90%  "

$ws.Cells.Item(15, 4).Value = "'2.671.96"
$ws.Cells.Item(15, 5).Value = "  +0.76%  "

$ws.Cells.Item(16, 4).Value = "'2.312.14"
$ws.Cells.Item(16, 5).Value = "  +1.45%  "

$ws.Cells.Item(17, 4).Value = "'0.788"
$ws.Cells.Item(17, 5).Value = "  +2.09%  "

$ws.Cells.Item(18, 4).Value = "'42.694.74"
$ws.Cells.Item(18, 5).Value = "  +0.24%  "

$ws.Cells.Item(19, 4).Value = "'12.21"
$ws.Cells.Item(19, 5).Value = "  -4.50%  "

$ws.Cells.Item(20, 4).Value = "'6.12"
$ws.Cells.Item(20, 5).Value = "  +2.27%  "

$ws.Cells.Item(21, 5).Value = "  -0.12%  "

$ws.Cells.Item(23, 5).Value = "  +6.77%  "

$ws.Cells.Item(24, 4).Value = "'235.04"

$ws.Cells.Item(25, 5).Value = "  -0.10%  "

$ws.Cells.Item(26, 5).Value = "  +1.23%  "

$ws.Cells.Item(27, 5).Value = "  -1.28%  "

$ws.Cells.Item(28, 5).Value = "  +15.09%  "

$ws.Cells.Item(29, 4).Value = "'166.32"
$ws.Cells.Item(29, 5).Value = "  -0.58%  "

$ws.Cells.Item(30, 4).Value = "'9.12"
$ws.Cells.Item(30, 5).Value = "  +1.66%  "

$ws.Cells.Item(31, 4).Value = "'32.09"
$ws.Cells.Item(31, 5).Value = "  -2.07%  "

$ws.Cells.Item(32, 5).Value = "  -0.01%  "

$ws.Cells.Item(33, 4).Value = "'4.99"
$ws.Cells.Item(33, 5).Value = "  +0.89%  "

$ws.Cells.Item(34, 4).Value = "'17.64"
$ws.Cells.Item(34, 5).Value = "  -0.54%  "

$ws.Cells.Item(35, 4).Value = "'4.46"
$ws.Cells.Item(35, 5).Value = "  +0.04%  "

$ws.Cells.Item(36, 4).Value = "'0.0697"
$ws.Cells.Item(36, 5).Value = "  +1.71%  "

$ws.Cells.Item(37, 5).Value = "  -0.91%  "

$ws.Cells.Item(38, 2).Value = "Kaspa"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(38, 4).Value = "'0.100"
$ws.Cells.Item(38, 5).Value = "  +0.14%  "

$ws.Cells.Item(39, 2).Value = "ARBITRUM"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(39, 4).Value = "'1.77"
$ws.Cells.Item(39, 5).Value = "  +2.10%  "

$ws.Cells.Item(40, 5).Value = "  +1.08%  "

$ws.Cells.Item(41, 5).Value = "  -0.48%  "

$ws.Cells.Item(42, 4).Value = "'20.76"
$ws.Cells.Item(42, 5).Value = "  +13.88%  "

$ws.Cells.Item(43, 4).Value = "'1.925.51"
$ws.Cells.Item(43, 5).Value = "  -3.32%  "

$ws.Cells.Item(44, 5).Value = "  -0.17%  "

$ws.Cells.Item(45, 4).Value = "'10.12"
$ws.Cells.Item(45, 5).Value = "  -0.46%  "

$ws.Cells.Item(46, 5).Value = "  -2.64%  "

$ws.Cells.Item(47, 4).Value = "'2.74"
$ws.Cells.Item(47, 5).Value = "  -0.68%  "

$ws.Cells.Item(48, 4).Value = "'2.88"
$ws.Cells.Item(48, 5).Value = "  +1.99%  "

$ws.Cells.Item(49, 4).Value = "'2.540.26"
$ws.Cells.Item(49, 5).Value = "  +0.86%  "

$ws.Cells.Item(50, 4).Value = "'53.26"
$ws.Cells.Item(50, 5).Value = "  -0.19%  "

$ws.Cells.Item(51, 4).Value = "'72.06"
$ws.Cells.Item(51, 5).Value = "  +1.90%  "
